# Update Oracle jobs data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: only the Apply_Link formula changes (job 31613 -> 31614)
$ws.Range("E3").Formula = '=HYPERLINK("https://estm.fa.em2.oraclecloud.com/hcmUI/CandidateExperience/en/sites/CX_1/job/31614/?location=India&locationId=300000000440677&locationLevel=country&mode=location", "Apply")'

# Row 4: title, posting date, and Apply_Link formula change
$ws.Range("B4").Value = "Finance Analyst"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "01/26/2026"
$ws.Range("E4").Formula = '=HYPERLINK("https://estm.fa.em2.oraclecloud.com/hcmUI/CandidateExperience/en/sites/CX_1/job/31468/?location=India&locationId=300000000440677&locationLevel=country&mode=location", "Apply")'

# Row 5: title and Apply_Link formula change (location/date stay the same)
$ws.Range("B5").Value = "Project Accounting & Financial Management Officer"
$ws.Range("E5").Formula = '=HYPERLINK("https://estm.fa.em2.oraclecloud.com/hcmUI/CandidateExperience/en/sites/CX_1/job/31469/?location=India&locationId=300000000440677&locationLevel=country&mode=location", "Apply")'

# Row 6 (old "Project Accounting & Financial Management Officer" entry) is removed entirely
$ws.Rows.Item(6).Delete()
